# Remove the standalone "9. " paragraph that sits right before the
# "Contenu de départ" heading. It is a leftover numbered-list remnant
# (pasted from an HTML/code block) whose text is "9. " followed by a
# single non-breaking space and nothing else.

$d = $word.ActiveDocument

$nbsp = [char]0x00A0
$cr = [char]13
$target = "{0}{1}{2}" -f "9. ", $nbsp, $cr

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq $target) {
        $p.Range.Delete()
        break
    }
}
